$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 104
$ws.Range("F4").Value = 31
$ws.Range("F5").Value = 8204
$ws.Range("F8").Value = 94
$ws.Range("F9").Value = 7144
$ws.Range("F10").Value = 1136
$ws.Range("F11").Value = 554
$ws.Range("F12").Value = 495
$ws.Range("F14").Value = 712
$ws.Range("F16").Value = 307
$ws.Range("F17").Value = 161
$ws.Range("F21").Value = 91
$ws.Range("F22").Value = 11721
$ws.Range("F24").Value = 132
$ws.Range("F25").Value = 2310
$ws.Range("F27").Value = 3227
$ws.Range("F29").Value = 2728
$ws.Range("F33").Value = 47
$ws.Range("F34").Value = 338
$ws.Range("F35").Value = 1629
$ws.Range("F36").Value = 75
$ws.Range("F37").Value = 108
$ws.Range("F38").Value = 5853
$ws.Range("F40").Value = 1801
$ws.Range("F41").Value = 1247
$ws.Range("F42").Value = 849
$ws.Range("F44").Value = 189
$ws.Range("F47").Value = 1535

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 26
$ws.Range("F8").Value = 253
$ws.Range("F20").Value = 68

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 254
$ws.Range("F3").Value = 396
$ws.Range("F4").Value = 2

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 254
$ws.Range("F4").Value = 396
$ws.Range("F7").Value = 8204
$ws.Range("F10").Value = 94
$ws.Range("F11").Value = 7144
$ws.Range("F12").Value = 7144
$ws.Range("F13").Value = 1136
$ws.Range("F14").Value = 554
$ws.Range("F15").Value = 495
$ws.Range("F16").Value = 712
$ws.Range("F18").Value = 307
$ws.Range("F19").Value = 161
$ws.Range("F21").Value = 253
$ws.Range("F22").Value = 91
$ws.Range("F25").Value = 11722
$ws.Range("F28").Value = 132
$ws.Range("F29").Value = 2310
$ws.Range("F30").Value = 2310
$ws.Range("F31").Value = 3227
$ws.Range("F32").Value = 2728
$ws.Range("F35").Value = 47
$ws.Range("F37").Value = 338
$ws.Range("F38").Value = 1629
$ws.Range("F39").Value = 75
$ws.Range("F40").Value = 108
$ws.Range("F41").Value = 5853
$ws.Range("F42").Value = 68
$ws.Range("F43").Value = 1801
$ws.Range("F45").Value = 1247
$ws.Range("F46").Value = 849
$ws.Range("F47").Value = 189
$ws.Range("F50").Value = 1535
